# Update "想去人数" (want-to-go count) values in column F for the sheets
# that hold the conference data: "展览" and "全部类型".
# Sheets "演出" and "本地生活" only contain header rows, so nothing to do there.

$wb = $excel.ActiveWorkbook

# Row number (in both sheets) -> new value for column F
$updates = @{
    2  = 11758
    3  = 11400
    4  = 607
    6  = 1029
    8  = 71
    9  = 46
    11 = 10787
    12 = 4172
    13 = 17
    14 = 8
    16 = 2471
    17 = 1054
    20 = 453
    21 = 11148
    22 = 10933
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
